$wb = $excel.ActiveWorkbook

# Insert a new column before column N (14) on the "Repayment Schedule" sheet
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# Select a cell and make this the active sheet/view
$wsSchedule.Activate()
$wsSchedule.Range("R12").Select()

# Restore the previous selection on the "Transactions" sheet (no longer the active tab)
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("D3").Select()

# Re-activate the Repayment Schedule sheet so it ends up as the selected/active tab
$wsSchedule.Activate()
